$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# AF2:AI2 currently hold the text "Missing" (inline/shared string cells).
# Replace them with the corrected numeric-looking text values, making sure
# Excel keeps them as TEXT (not auto-converted to numbers) - this mirrors
# the other text-typed numeric cells already present in the row and avoids
# the "leading zero" / auto-numeric-conversion bug the commit fixes.

$cells = @("AF2", "AG2", "AH2", "AI2")
$values = @("0.029592333", "8.69794377", "70.5995828962491", "18.5")

for ($i = 0; $i -lt $cells.Count; $i++) {
    $rng = $ws.Range($cells[$i])
    $rng.NumberFormat = "@"
    $rng.Value = $values[$i]
    $rng.Style = "Normal"
}
